$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.936.01"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").Value = "3.445.00"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.28%  "

$ws.Range("D9").Value = "3.445.82"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -2.00%  "

$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("E12").Value = "  -2.70%  "

$ws.Range("D13").Value = "4.040.88"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.50%  "

$ws.Range("D16").Value = "65.917.98"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "3.443.58"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.529"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("E26").Value = "  +3.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.61%  "

$ws.Range("E28").Value = "  +3.49%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.04%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.69%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").Value = "2.765.14"
$ws.Range("E42").Value = "  +2.97%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0680"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0290"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "324.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
